$d = $word.ActiveDocument

# Grab the full OOXML package for the document body/content.
$xml = $d.Content.WordOpenXML

# The edit simply removes every <w:contextualSpacing w:val="0"/> element
# (there is exactly one per <w:pPr> block in this document, always with
# the same w:val="0" value) without touching anything else.
$needle = '<w:contextualSpacing w:val="0"/>'
$newXml = $xml.Replace($needle, "")

# Push the modified package XML back into the document, replacing the
# whole content in one shot.
$d.Content.InsertXML($newXml)
